$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.992.42'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.631.59'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.990'
$ws.Range("E4").Value = '  -1.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.69'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.992'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.256'
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0631'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.73'
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.24'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '1.863.30'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '1.622.17'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.551'
$ws.Range("E15").Value = '  -1.78%  '
$ws.Range("D16").Value = '0.0₃0762'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.14'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '25.993.84'
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.991'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.45'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.91'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.98'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.34'
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.990'
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.40'
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.124'
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.88'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.58'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0496'
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.38'
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.903'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").Value = '1.136.81'
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.547'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.989'
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.99'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.785'
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("D45").Value = '1.773.50'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").Value = '0.0₆0107'
$ws.Range("E46").Value = '  -4.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.76'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("E49").Value = '  +4.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.415'
$ws.Range("E51").Value = '  -0.29%  '
